# Generate Report for Handoff
# Replaces the old job/file identifiers (uuid "1fc347ef-..." + xliff hash
# "5e2f077c...") with the new ones ("ab1ab29d-..." / "f9fd2827...") across
# the Overview / zh-cn / de-de sheets, and bumps the associated "Latest
# Handoff" timestamps forward (handoff generation re-run).

$wb = $excel.ActiveWorkbook

$oldGuid = "1fc347ef-60fc-4297-95ec-5dd81550e5c4"
$newGuid = "ab1ab29d-8a17-450c-9624-2b30c65f8987"
$oldHash = "5e2f077cfdcdc2a5f6110e77d0caebd9351c923d"
$newHash = "f9fd2827ba778bd26ebbe2cb5e0181c794c1f0e5"

$newMdName  = "$newGuid.md"
$newMdPath  = "e2e\$newGuid.md"
$newDate1   = "2016-08-30 07:01:21"
$newZhXlf   = "$newGuid.$newHash.zh-cn.xlf"
$newDeXlf   = "$newGuid.$newHash.de-de.xlf"
$newZhDate  = "2016-08-30 07:01:16"

# ---- Overview sheet ----
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("A2").Value = $newMdName
$wsOverview.Range("B2").Value = $newMdPath
$wsOverview.Range("G2").Value = $newDate1
# NB: update the hyperlink display text via the foreach enumerator -
# Hyperlinks.Item(n) returns a detached instance whose property writes
# are silently dropped (and re-reading it back yields empty strings);
# iterating the collection gives the live object instead.
foreach ($hl in $wsOverview.Hyperlinks) {
    $hl.TextToDisplay = $newMdPath
}

# ---- zh-cn sheet ----
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("A2").Value = $newMdName
$wsZhCn.Range("G2").Value = $newZhXlf
$wsZhCn.Range("H2").Value = $newZhDate
foreach ($hl in $wsZhCn.Hyperlinks) {
    $hl.TextToDisplay = $newMdName
}

# ---- de-de sheet ----
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("A2").Value = $newMdName
$wsDeDe.Range("G2").Value = $newDeXlf
$wsDeDe.Range("H2").Value = $newDate1
foreach ($hl in $wsDeDe.Hyperlinks) {
    $hl.TextToDisplay = $newMdName
}
